$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data set (interleaved by period, ascending 1907..2001, PEDRO then LUIS for
# each period). PEDRO's "Salario Basico" (column G) is updated to 2500000 for all
# his rows; LUIS's values are unchanged aside from the row position.
$rows = @(
    @{ Row = 16; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "1907"; ValorMora = 100000; Salario = 2500000 },
    @{ Row = 17; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "1907"; ValorMora = 33125;  Salario = 828116  },
    @{ Row = 18; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "1908"; ValorMora = 100000; Salario = 2500000 },
    @{ Row = 19; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "1908"; ValorMora = 33125;  Salario = 828116  },
    @{ Row = 20; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "1909"; ValorMora = 100000; Salario = 2500000 },
    @{ Row = 21; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "1909"; ValorMora = 33125;  Salario = 828116  },
    @{ Row = 22; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "1910"; ValorMora = 100000; Salario = 2500000 },
    @{ Row = 23; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "1910"; ValorMora = 33125;  Salario = 828116  },
    @{ Row = 24; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "1911"; ValorMora = 100000; Salario = 2500000 },
    @{ Row = 25; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "1911"; ValorMora = 33125;  Salario = 828116  },
    @{ Row = 26; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "1912"; ValorMora = 100000; Salario = 2500000 },
    @{ Row = 27; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "1912"; ValorMora = 33125;  Salario = 828116  },
    @{ Row = 28; TipoDoc = "CC"; NumDoc = "73148203";   Nombre = "PEDRO UTRIA MONSALVE";          Periodo = "2001"; ValorMora = 70000;  Salario = 2500000 },
    @{ Row = 29; TipoDoc = "CC"; NumDoc = "1152226681"; Nombre = "LUIS MIGUELL UTRIA MORALESS";   Periodo = "2001"; ValorMora = 23187;  Salario = 828116  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.TipoDoc
    $ws.Range("C$n").Value = $r.NumDoc
    $ws.Range("D$n").Value = $r.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.ValorMora
    $ws.Range("G$n").Value = $r.Salario
}
